# Testcase 12 (12.1, 12.2): split the old "BillingCreateOrders" sheet into
# two sheets - a new "BillCreateCustomers" sheet (customer/billing-account
# setup fields) feeding a trimmed-down "BillCreateOrders" sheet (order
# creation fields only, duplicate trailing block removed).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new "BillCreateCustomers" sheet right before the existing
#    "BillingCreateOrders" sheet.
# ---------------------------------------------------------------------
$orders = $wb.Worksheets.Item("BillingCreateOrders")
$new = $wb.Worksheets.Add($orders)
$new.Name = "BillCreateCustomers"

# Sheet handles returned before a worksheet insert/rename can go stale
# (they track position, not identity), so re-resolve by name before use.
$custs = $wb.Worksheets.Item("BillCreateCustomers")

$custs.Range("A1").Value = "admin"
$custs.Range("B1").Value = "Webdata@123"
$custs.Range("C1").Value = "Web Data US"
$custs.Range("D1").Value = "Web Data US"
$custs.Range("E1").Value = "Direct Customer"
$custs.Range("F1").Value = "Billing Customer1"
$custs.Range("G1").Value = "Monthly"
$custs.Range("H1").Value = 1
$custs.Range("I1").Value = 1
$custs.Range("J1").Value = "Billing1"

$custs.Range("K1").NumberFormat = "0"
$custs.Range("K1").Value = 4111111111111150

$custs.Range("L1").NumberFormat = "mmm-yy"
$custs.Range("L1").Value = 43862

$custs.Range("M1").NumberFormat = "dd/mm/yyyy"
$custs.Range("M1").Value = 36892

$custs.Range("N1").Value = "Billing Customer2"
$custs.Range("O1").Value = "Billing2"
$custs.Range("P1").Value = 15

$custs.Range("Q1").NumberFormat = "dd/mm/yyyy"
$custs.Range("Q1").Value = 36892

# Leftover formatted-but-empty date cells (same pattern Excel leaves behind
# after a partial column copy/paste of the date columns).
foreach ($addr in @("U1", "Y1", "AC1", "AD1", "AH1", "AI1", "AM1", "AN1")) {
    $custs.Range($addr).NumberFormat = "dd/mm/yyyy"
}

$custs.Columns.Item(11).ColumnWidth = 19.83203125
$custs.Columns.Item(12).ColumnWidth = 19.83203125

$custs.Range("P1").Select()

# ---------------------------------------------------------------------
# 2) Rename "BillingCreateOrders" -> "BillCreateOrders" and rewrite its
#    row so it only keeps the order-creation fields (the customer/account
#    fields now live on BillCreateCustomers, and the duplicated trailing
#    "Monthly/post paid" block collapses to a single occurrence).
# ---------------------------------------------------------------------
$orders = $wb.Worksheets.Item("BillingCreateOrders")
$orders.Name = "BillCreateOrders"

$orders = $wb.Worksheets.Item("BillCreateOrders")
$orders.Rows.Item(1).Clear()

$orders.Range("A1").Value = "admin"
$orders.Range("B1").Value = "Webdata@123"
$orders.Range("C1").Value = "Web Data US"
$orders.Range("D1").Value = "Billing Category"
$orders.Range("E1").Value = "Billing Flat"
$orders.Range("F1").Value = "BF01"

$orders.Range("G1").NumberFormat = "dd/mm/yyyy"
$orders.Range("G1").Value = 36892

$orders.Range("H1").Value = 20
$orders.Range("I1").Value = "Billing Graduated"
$orders.Range("J1").Value = "BG01"

$orders.Range("K1").NumberFormat = "dd/mm/yyyy"
$orders.Range("K1").Value = 36892

$orders.Range("L1").Value = 5
$orders.Range("M1").Value = "Billing Customer2"
$orders.Range("N1").Value = "One Time"
$orders.Range("O1").Value = "post paid"

$orders.Range("P1").NumberFormat = "dd/mm/yyyy"
$orders.Range("P1").Value = 36892
$orders.Range("Q1").NumberFormat = "dd/mm/yyyy"
$orders.Range("Q1").Value = 36892

$orders.Range("R1").Value = 14
$orders.Range("S1").Value = "Monthly"
$orders.Range("T1").Value = "pre paid"

$orders.Range("U1").NumberFormat = "dd/mm/yyyy"
$orders.Range("U1").Value = 36892
$orders.Range("V1").NumberFormat = "dd/mm/yyyy"
$orders.Range("V1").Value = 36892

$orders.Range("W1").Value = 3
$orders.Range("X1").Value = "Billing Customer1"
$orders.Range("Y1").Value = "Monthly"
$orders.Range("Z1").Value = "post paid"

$orders.Range("AA1").NumberFormat = "dd/mm/yyyy"
$orders.Range("AA1").Value = 36892
$orders.Range("AB1").NumberFormat = "dd/mm/yyyy"
$orders.Range("AB1").Value = 36892

$orders.Range("AC1").Value = 14
$orders.Range("AD1").Value = 2

# BillCreateOrders becomes the active tab/sheet.
$orders = $wb.Worksheets.Item("BillCreateOrders")
$orders.Activate()
$orders.Range("E1").Select()
